$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.289955139160156
$ws.Range("B1").Value = 2.527620792388916
$ws.Range("C1").Value = 2.610470771789551
$ws.Range("D1").Value = 3.329218626022339
$ws.Range("E1").Value = 2.45249605178833
